# conductor_diagnostic.xlsx - rename the two worksheets to their final,
# more descriptive names. Excel automatically rewrites every formula that
# referenced the old sheet names (e.g. "Space!A1" -> "Spatial_distribution!A1").
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Space").Name = "Spatial_distribution"
$wb.Worksheets.Item("Time").Name = "Time_evolution"
